$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Roll the daily price table forward: previous "today" (46004) becomes "yesterday" (46004->46004... )
# Step 1: copy the CURRENT "today" row values (date 46004) into the corresponding "yesterday" row (date 46003 -> 46004),
#         using the values that were present before this update (i.e. yesterday row gets what used to be today).
$ws.Range("A11").Value = 46004
$ws.Range("D11").Value = 159.78
$ws.Range("E11").Value = 158.4
$ws.Range("F11").Value = 168.4
$ws.Range("G11").Value = 158.51

$ws.Range("A12").Value = 46004
$ws.Range("D12").Value = 159.78
$ws.Range("E12").Value = 158.4
$ws.Range("F12").Value = 168.4
$ws.Range("G12").Value = 158.51

$ws.Range("A13").Value = 46004
$ws.Range("D13").Value = 161.44
$ws.Range("E13").Value = 161.12
$ws.Range("F13").Value = 171.12
$ws.Range("G13").Value = 161.63999999999999

$ws.Range("A18").Value = 46004
$ws.Range("D18").Value = 164.39
$ws.Range("E18").Value = 164.48
$ws.Range("F18").Value = 174.48

$ws.Range("A27").Value = 46004
$ws.Range("D27").Value = 160.37
$ws.Range("E27").Value = 160.22
$ws.Range("F27").Value = 169.82
$ws.Range("G27").Value = 161.38999999999999

$ws.Range("A28").Value = 46004
$ws.Range("D28").Value = 166.47
$ws.Range("E28").Value = 165.23
$ws.Range("F28").Value = 175.23

$ws.Range("A29").Value = 46004
$ws.Range("D29").Value = 166.24
$ws.Range("E29").Value = 165.55
$ws.Range("F29").Value = 175.55

$ws.Range("A30").Value = 46004
$ws.Range("D30").Value = 166.85
$ws.Range("E30").Value = 164.98
$ws.Range("F30").Value = 174.98
$ws.Range("G30").Value = 164.75

$ws.Range("A31").Value = 46004
$ws.Range("D31").Value = 165.65
$ws.Range("E31").Value = 166.41
$ws.Range("F31").Value = 176.41

$ws.Range("A36").Value = 46004
$ws.Range("D36").Value = 159.84
$ws.Range("E36").Value = 158.03
$ws.Range("F36").Value = 167.04

$ws.Range("A42").Value = 46004
$ws.Range("D42").Value = 165.72
$ws.Range("E42").Value = 165.12
$ws.Range("F42").Value = 175.12

$ws.Range("A43").Value = 46004
$ws.Range("D43").Value = 165.43
$ws.Range("E43").Value = 165.54
$ws.Range("F43").Value = 175.54

$ws.Range("A49").Value = 46004
$ws.Range("D49").Value = 160.74
$ws.Range("E49").Value = 159.80000000000001
$ws.Range("F49").Value = 169.8

$ws.Range("A50").Value = 46004
$ws.Range("D50").Value = 160.55000000000001
$ws.Range("E50").Value = 159.88999999999999
$ws.Range("F50").Value = 169.89

$ws.Range("A60").Value = 46004
$ws.Range("D60").Value = 175.23
$ws.Range("E60").Value = 175.76
$ws.Range("F60").Value = 185.76

$ws.Range("A61").Value = 46004
$ws.Range("D61").Value = 163.41
$ws.Range("E61").Value = 170.91
$ws.Range("F61").Value = 180.91

$ws.Range("A62").Value = 46004
$ws.Range("D62").Value = 166

$ws.Range("A63").Value = 46004
$ws.Range("D63").Value = 165.01
$ws.Range("E63").Value = 165.18

$ws.Range("A64").Value = 46004
$ws.Range("D64").Value = 160.91999999999999
$ws.Range("E64").Value = 161.22999999999999
$ws.Range("F64").Value = 171.23

$ws.Range("A65").Value = 46004
$ws.Range("D65").Value = 167.31
$ws.Range("E65").Value = 173.03

# Step 2: write the brand-new "today" values (date 46007) into the former "today" rows.
$ws.Range("A8").Value = 46007
$ws.Range("D8").Value = 159.01
$ws.Range("E8").Value = 158.24
$ws.Range("F8").Value = 168.24
$ws.Range("G8").Value = 158.36000000000001

$ws.Range("A9").Value = 46007
$ws.Range("D9").Value = 159.01
$ws.Range("E9").Value = 158.24
$ws.Range("F9").Value = 168.24
$ws.Range("G9").Value = 158.36000000000001

$ws.Range("A10").Value = 46007
$ws.Range("D10").Value = 160.69
$ws.Range("E10").Value = 160.87
$ws.Range("F10").Value = 170.87
$ws.Range("G10").Value = 161.38999999999999

$ws.Range("A17").Value = 46007
$ws.Range("D17").Value = 163.63999999999999
$ws.Range("E17").Value = 164.17
$ws.Range("F17").Value = 174.17

$ws.Range("A22").Value = 46007
$ws.Range("D22").Value = 159.49
$ws.Range("E22").Value = 159.97999999999999
$ws.Range("F22").Value = 169.57
$ws.Range("G22").Value = 161.13999999999999

$ws.Range("A23").Value = 46007
$ws.Range("D23").Value = 165.72
$ws.Range("E23").Value = 164.98
$ws.Range("F23").Value = 174.98

$ws.Range("A24").Value = 46007
$ws.Range("D24").Value = 165.49
$ws.Range("E24").Value = 165.3
$ws.Range("F24").Value = 175.3

$ws.Range("A25").Value = 46007
$ws.Range("D25").Value = 166.1
$ws.Range("E25").Value = 164.73
$ws.Range("F25").Value = 174.73
$ws.Range("G25").Value = 164.5

$ws.Range("A26").Value = 46007
$ws.Range("D26").Value = 164.9
$ws.Range("E26").Value = 166.16
$ws.Range("F26").Value = 176.16

$ws.Range("A35").Value = 46007
$ws.Range("D35").Value = 158.87
$ws.Range("E35").Value = 157.79
$ws.Range("F35").Value = 166.79

$ws.Range("A40").Value = 46007
$ws.Range("D40").Value = 164.97
$ws.Range("E40").Value = 164.82
$ws.Range("F40").Value = 174.82

$ws.Range("A41").Value = 46007
$ws.Range("D41").Value = 164.69
$ws.Range("E41").Value = 165.24
$ws.Range("F41").Value = 175.24

$ws.Range("A47").Value = 46007
$ws.Range("D47").Value = 160.15
$ws.Range("E47").Value = 159.4
$ws.Range("F47").Value = 169.4

$ws.Range("A48").Value = 46007
$ws.Range("D48").Value = 159.96
$ws.Range("E48").Value = 159.49
$ws.Range("F48").Value = 169.49

$ws.Range("A54").Value = 46007
$ws.Range("D54").Value = 174.47
$ws.Range("E54").Value = 175.52
$ws.Range("F54").Value = 185.52

$ws.Range("A55").Value = 46007
$ws.Range("D55").Value = 162.66
$ws.Range("E55").Value = 170.66
$ws.Range("F55").Value = 180.66

$ws.Range("A56").Value = 46007
$ws.Range("D56").Value = 165.03

$ws.Range("A57").Value = 46007
$ws.Range("D57").Value = 164.04
$ws.Range("E57").Value = 164.93

$ws.Range("A58").Value = 46007
$ws.Range("D58").Value = 159.94999999999999
$ws.Range("E58").Value = 160.97999999999999
$ws.Range("F58").Value = 170.98

$ws.Range("A59").Value = 46007
$ws.Range("D59").Value = 166.56
$ws.Range("E59").Value = 172.78
